$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1690.8536
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 1912.1428
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 5736.428400000001
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -6072.428400000001

# Row 111
$ws.Range("H111").Value = 1750
$ws.Range("I111").Value = 2050
$ws.Range("J111").Value = 1600
$ws.Range("K111").Value = 6150
$ws.Range("L111").Value = 4800
$ws.Range("M111").Value = -3083
$ws.Range("N111").Value = -10934

# Row 113
$ws.Range("H113").Value = 287001.25
$ws.Range("I113").Value = 287001.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 287001.25
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -283747.25

# Row 116
$ws.Range("H116").Value = 5767398
$ws.Range("I116").Value = 6018024
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 6018024
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -6014582
$ws.Range("N116").Value = -9884

# Row 132
$ws.Range("H132").Value = 276592.06
$ws.Range("I132").Value = 338395.5
$ws.Range("J132").Value = 54099.7
$ws.Range("K132").Value = 1015186.5
$ws.Range("L132").Value = 162299.1
$ws.Range("M132").Value = -1012656.5
$ws.Range("N132").Value = -167359.1


$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 12500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 12500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 12500
$ws.Range("N24").Value = -13248

# Row 32
$ws.Range("H32").Value = 16870.367
$ws.Range("I32").Value = 1963.6094
$ws.Range("J32").Value = 255378.5
$ws.Range("K32").Value = 1963.6094
$ws.Range("L32").Value = 255378.5
$ws.Range("M32").Value = -1676.6094
$ws.Range("N32").Value = -255952.5

# Row 100
$ws.Range("H100").Value = 12500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 12500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 12500
$ws.Range("N100").Value = -14664

# Row 122
$ws.Range("H122").Value = 1502.1522
$ws.Range("I122").Value = 1292.3793
$ws.Range("J122").Value = 1860
$ws.Range("K122").Value = 3877.1379
$ws.Range("L122").Value = 5580
$ws.Range("M122").Value = -1427.1379
$ws.Range("N122").Value = -10480


$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 2649.8333
$ws.Range("I12").Value = 999.6667
$ws.Range("J12").Value = 4300
$ws.Range("K12").Value = 999.6667
$ws.Range("L12").Value = 4300
$ws.Range("M12").Value = -831.6667
$ws.Range("N12").Value = -4636

# Row 75
$ws.Range("H75").Value = 98794.89999999999
$ws.Range("I75").Value = 5542.6
$ws.Range("J75").Value = 192047.2
$ws.Range("K75").Value = 5542.6
$ws.Range("L75").Value = 192047.2
$ws.Range("M75").Value = -4606.6
$ws.Range("N75").Value = -193919.2

# Row 78
$ws.Range("H78").Value = 98794.89999999999
$ws.Range("I78").Value = 5542.6
$ws.Range("J78").Value = 192047.2
$ws.Range("K78").Value = 16627.8
$ws.Range("L78").Value = 576141.6000000001
$ws.Range("M78").Value = -11947.8
$ws.Range("N78").Value = -585501.6000000001

# Row 94
$ws.Range("H94").Value = 1265.6666
$ws.Range("I94").Value = 1511.875
$ws.Range("J94").Value = 773.25
$ws.Range("K94").Value = 1511.875
$ws.Range("L94").Value = 773.25
$ws.Range("M94").Value = -1060.875
$ws.Range("N94").Value = -1675.25

# Row 107
$ws.Range("H107").Value = 804.7143
$ws.Range("I107").Value = 690
$ws.Range("J107").Value = 1493
$ws.Range("K107").Value = 690
$ws.Range("L107").Value = 1493
$ws.Range("M107").Value = 1230
$ws.Range("N107").Value = -5333


$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 50.5
$ws.Range("I6").Value = 50.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 50.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 62.5

# Row 16
$ws.Range("H16").Value = 1499.6666
$ws.Range("I16").Value = 1449.5
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 1449.5
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -1162.5
$ws.Range("N16").Value = -2174

# Row 31
$ws.Range("H31").Value = 1387.88
$ws.Range("I31").Value = 1031.6818
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1031.6818
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -736.6818000000001
$ws.Range("N31").Value = -4590

# Row 34
$ws.Range("H34").Value = 1387.88
$ws.Range("I34").Value = 1031.6818
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1031.6818
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -829.6818000000001
$ws.Range("N34").Value = -4404

# Row 105
$ws.Range("H105").Value = 832
$ws.Range("I105").Value = 897
$ws.Range("J105").Value = 799.5
$ws.Range("K105").Value = 897
$ws.Range("L105").Value = 799.5
$ws.Range("M105").Value = 850
$ws.Range("N105").Value = -4293.5

# Row 113
$ws.Range("H113").Value = 1499.6666
$ws.Range("I113").Value = 1449.5
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 1449.5
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 720.5
$ws.Range("N113").Value = -5940


$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1353.8695
$ws.Range("I5").Value = 702
$ws.Range("J5").Value = 2065
$ws.Range("K5").Value = 2106
$ws.Range("L5").Value = 6195
$ws.Range("M5").Value = -1994
$ws.Range("N5").Value = -6419

# Row 64
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 3666.6667
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 11000.0001
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -11540.0001

# Row 67
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 3666.6667
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 11000.0001
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -12872.0001

# Row 114
$ws.Range("H114").Value = 1146.3077
$ws.Range("I114").Value = 384.6
$ws.Range("J114").Value = 1622.375
$ws.Range("K114").Value = 1153.8
$ws.Range("L114").Value = 4867.125
$ws.Range("M114").Value = 2100.2
$ws.Range("N114").Value = -11375.125

# Row 122
$ws.Range("H122").Value = 1119.4
$ws.Range("I122").Value = 502
$ws.Range("J122").Value = 1531
$ws.Range("K122").Value = 4518
$ws.Range("L122").Value = 13779
$ws.Range("M122").Value = -2068
$ws.Range("N122").Value = -18679

# Row 135
$ws.Range("H135").Value = 1353.8695
$ws.Range("I135").Value = 702
$ws.Range("J135").Value = 2065
$ws.Range("K135").Value = 6318
$ws.Range("L135").Value = 18585
$ws.Range("M135").Value = -3783
$ws.Range("N135").Value = -23655


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2219.6667
$ws.Range("I102").Value = 1913.6
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 1913.6
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -291.5999999999999
$ws.Range("N102").Value = -6994

# Row 126
$ws.Range("H126").Value = 2386.0344
$ws.Range("I126").Value = 1644
$ws.Range("J126").Value = 2839.5
$ws.Range("K126").Value = 4932
$ws.Range("L126").Value = 8518.5
$ws.Range("M126").Value = -2462
$ws.Range("N126").Value = -13458.5


$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 50000
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -49887

# Row 28
$ws.Range("H28").Value = 50000
$ws.Range("I28").Value = 50000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 50000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -49768

# Row 37
$ws.Range("H37").Value = 50000
$ws.Range("I37").Value = 50000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 50000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -49893

# Row 40
$ws.Range("H40").Value = 2647.4482
$ws.Range("I40").Value = 1576.5
$ws.Range("J40").Value = 3055.4285
$ws.Range("K40").Value = 1576.5
$ws.Range("L40").Value = 3055.4285
$ws.Range("M40").Value = -1440.5
$ws.Range("N40").Value = -3327.4285

# Row 61
$ws.Range("H61").Value = 3619.1785
$ws.Range("I61").Value = 3534.2273
$ws.Range("J61").Value = 3930.6667
$ws.Range("K61").Value = 3534.2273
$ws.Range("L61").Value = 3930.6667
$ws.Range("M61").Value = -3332.2273
$ws.Range("N61").Value = -4334.6667

# Row 113
$ws.Range("H113").Value = 3619.1785
$ws.Range("I113").Value = 3534.2273
$ws.Range("J113").Value = 3930.6667
$ws.Range("K113").Value = 3534.2273
$ws.Range("L113").Value = 3930.6667
$ws.Range("M113").Value = -1364.2273
$ws.Range("N113").Value = -8270.6667

# Row 122
$ws.Range("H122").Value = 3025.0278
$ws.Range("I122").Value = 1900
$ws.Range("J122").Value = 3587.5417
$ws.Range("K122").Value = 5700
$ws.Range("L122").Value = 10762.6251
$ws.Range("M122").Value = -3250
$ws.Range("N122").Value = -15662.6251

# Row 136
$ws.Range("H136").Value = 4097.457
$ws.Range("I136").Value = 1832.16
$ws.Range("J136").Value = 9760.700000000001
$ws.Range("K136").Value = 5496.48
$ws.Range("L136").Value = 29282.1
$ws.Range("M136").Value = -2946.48
$ws.Range("N136").Value = -34382.10000000001


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 842.25
$ws.Range("I113").Value = 711.9231
$ws.Range("J113").Value = 1407
$ws.Range("K113").Value = 2135.7693
$ws.Range("L113").Value = 4221
$ws.Range("M113").Value = 34.23070000000007
$ws.Range("N113").Value = -8561

# Row 123
$ws.Range("H123").Value = 25000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 25000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

